$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.142.86"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "2.262.87"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").Value = "'495.70"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'128.53"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("D12").Value = "'4.86"
$ws.Range("E12").Value = "  +4.59%  "
$ws.Range("E13").Value = "  +5.25%  "
$ws.Range("D14").Value = "2.662.39"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").Value = "54.132.42"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "2.263.65"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D19").Value = "'4.13"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").Value = "'302.96"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'60.56"
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("D27").Value = "'172.22"
$ws.Range("E27").Value = "  +2.50%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "'5.95"
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("D30").Value = "0.0₃0689"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'0.942"
$ws.Range("E35").Value = "  +3.30%  "
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "'0.374"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "'1.40"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").Value = "'4.79"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").Value = "'124.39"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("D45").Value = "'0.543"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "'240.85"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "'0.934"
